$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1) Border styling for the "header row" spacer cells.
#    C1/F1 get a top+bottom thin border (no left/right).
#    D1/G1 get a top+bottom+right thin border (no left).
#    These combinations already exist in the workbook's border table
#    (index 4 and index 5 respectively), so we only need to create two
#    new cellXfs entries that reference them - achieved by clearing the
#    cell's borders first and then turning on just the edges we need.
# ---------------------------------------------------------------------

$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = -4142          # xlLineStyleNone - clear all edges
$c1.Borders.Item(8).LineStyle = 1      # xlEdgeTop - xlContinuous
$c1.Borders.Item(9).LineStyle = 1      # xlEdgeBottom - xlContinuous

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$c1.Copy()
$d1.PasteSpecial(-4122)                # xlPasteFormats - reuse the top+bottom style
$d1.Borders.Item(10).LineStyle = 1     # xlEdgeRight - xlContinuous

# Propagate the two established styles to sheet2's matching cells purely
# via format copy/paste so no extra intermediate styles get materialised.
$c1b = $ws2.Range("C1")
$c1b.Style = "Normal"
$c1.Copy()
$c1b.PasteSpecial(-4122)

$d1b = $ws2.Range("D1")
$d1b.Style = "Normal"
$d1.Copy()
$d1b.PasteSpecial(-4122)

$f1b = $ws2.Range("F1")
$f1b.Style = "Normal"
$c1.Copy()
$f1b.PasteSpecial(-4122)

$g1b = $ws2.Range("G1")
$g1b.Style = "Normal"
$d1.Copy()
$g1b.PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Anonymize "fedcore" -> "approach" in both sheets' header rows.
# ---------------------------------------------------------------------
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---------------------------------------------------------------------
# 3) Drop the stray empty inline-string cell G5 on computational_comparison.
# ---------------------------------------------------------------------
$ws2.Range("G5").ClearContents()
